$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDTPTUMCF")

$ws.Range("B19").Formula = "=B2"
$ws.Range("B20").Formula = "=B4"
$ws.Range("B21").Formula = "=B10"
$ws.Range("B22").Formula = "=B14"
$ws.Range("B23").Formula = "=B5"

$rng = $ws.Range("B19:B23")
$rng.Interior.Color = 14277081
